$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.085.88"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "1.972.86"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").Value = "329.69"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "0.4975"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").Value = "0.4188"
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "54.31"
$ws.Range("E9").Value = "  +3.78%  "
$ws.Range("D10").Value = "0.09281"
$ws.Range("E10").Value = "  +4.12%  "
$ws.Range("D11").Value = "1.094"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "22.75"
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("D13").Value = "1.988.89"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "7.876"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "6.450"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "91.62"
$ws.Range("E18").Value = "  -4.92%  "
$ws.Range("D19").Value = "0.06737"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").Value = "19.12"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").Value = "1.014"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "5.958"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "29.097.90"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "2.272"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "2.214.68"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "20.75"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "156.79"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").Value = "6.239"
$ws.Range("E29").Value = "  -4.62%  "
$ws.Range("D30").Value = "2.260"
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("D31").Value = "127.24"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "1.043"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("D33").Value = "0.09826"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").Value = "1.500"
$ws.Range("E34").Value = "  -4.16%  "
$ws.Range("D35").Value = "5.798"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("D36").Value = "3.750"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").Value = "0.02418"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").Value = "1.322"
$ws.Range("E38").Value = "  +2.67%  "
$ws.Range("D39").Value = "0.06395"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "9.037"
$ws.Range("E40").Value = "  -5.83%  "
$ws.Range("D41").Value = "0.6464"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").Value = "11.45"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("D43").Value = "0.2000"
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("D44").Value = "1.014"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "0.6181"
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("D46").Value = "1.349"
$ws.Range("E46").Value = "  +6.10%  "
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "2.172"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "0.06958"
$ws.Range("E51").Value = "  -0.62%  "
